$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1157.4
$ws.Range("I18").Value = 1157.4
$ws.Range("K18").Value = 1157.4
$ws.Range("M18").Value = -873.4000000000001

$ws.Range("H41").Value = 937.375
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 1083.1666
$ws.Range("K41").Value = 500
$ws.Range("L41").Value = 1083.1666
$ws.Range("M41").Value = -60
$ws.Range("N41").Value = -1963.1666

$ws.Range("H80").Value = 1840.3125
$ws.Range("I80").Value = 1449.6364
$ws.Range("J80").Value = 2699.8
$ws.Range("K80").Value = 4348.9092
$ws.Range("L80").Value = 8099.400000000001
$ws.Range("M80").Value = -3350.9092
$ws.Range("N80").Value = -10095.4

$ws.Range("H83").Value = 1840.3125
$ws.Range("I83").Value = 1449.6364
$ws.Range("J83").Value = 2699.8
$ws.Range("K83").Value = 13046.7276
$ws.Range("L83").Value = 24298.2
$ws.Range("M83").Value = -8054.7276
$ws.Range("N83").Value = -34282.2

$ws.Range("H98").Value = 1040.0555
$ws.Range("I98").Value = 948.2941
$ws.Range("K98").Value = 948.2941
$ws.Range("M98").Value = 549.7059

$ws.Range("H106").Value = 9868.9
$ws.Range("I106").Value = 4781.5
$ws.Range("K106").Value = 4781.5
$ws.Range("M106").Value = -4150.5

$ws.Range("H122").Value = 1040.0555
$ws.Range("I122").Value = 948.2941
$ws.Range("K122").Value = 2844.8823
$ws.Range("M122").Value = -394.8822999999998

$ws.Range("H132").Value = 18686.309
$ws.Range("I132").Value = 18686.309
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 56058.927
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -53528.927

$ws.Range("H135").Value = 7556.3823
$ws.Range("I135").Value = 6478.407
$ws.Range("J135").Value = 11714.286
$ws.Range("K135").Value = 58305.663
$ws.Range("L135").Value = 105428.574
$ws.Range("M135").Value = -55770.663
$ws.Range("N135").Value = -110498.574

$ws.Range("H138").Value = 6353.877
$ws.Range("J138").Value = 6140.375
$ws.Range("L138").Value = 18421.125
$ws.Range("N138").Value = -28701.125

$ws.Range("H141").Value = 4433.6206
$ws.Range("I141").Value = 1612.591
$ws.Range("J141").Value = 13299.714
$ws.Range("K141").Value = 4837.772999999999
$ws.Range("L141").Value = 39899.142
$ws.Range("M141").Value = 342.2270000000008
$ws.Range("N141").Value = -50259.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2307.3809
$ws.Range("I45").Value = 2153.1667
$ws.Range("K45").Value = 2153.1667
$ws.Range("M45").Value = -1776.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1332.5
$ws.Range("I107").Value = 1201
$ws.Range("J107").Value = 1990
$ws.Range("K107").Value = 1201
$ws.Range("L107").Value = 1990
$ws.Range("M107").Value = 719
$ws.Range("N107").Value = -5830

$ws.Range("H134").Value = 6434257.5
$ws.Range("I134").Value = 6019623.5
$ws.Range("K134").Value = 18058870.5
$ws.Range("M134").Value = -18056335.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71519.57000000001
$ws.Range("I7").Value = 100080.1
$ws.Range("K7").Value = 100080.1
$ws.Range("M7").Value = -99967.10000000001

$ws.Range("H31").Value = 345014.94
$ws.Range("I31").Value = 629797.6
$ws.Range("J31").Value = 3275.7144
$ws.Range("K31").Value = 629797.6
$ws.Range("L31").Value = 3275.7144
$ws.Range("M31").Value = -629502.6
$ws.Range("N31").Value = -3865.7144

$ws.Range("H34").Value = 345014.94
$ws.Range("I34").Value = 629797.6
$ws.Range("J34").Value = 3275.7144
$ws.Range("K34").Value = 629797.6
$ws.Range("L34").Value = 3275.7144
$ws.Range("M34").Value = -629595.6
$ws.Range("N34").Value = -3679.7144

$ws.Range("H58").Value = 6007796.5
$ws.Range("I58").Value = 23814222
$ws.Range("K58").Value = 23814222
$ws.Range("M58").Value = -23814019

$ws.Range("H132").Value = 2947
$ws.Range("I132").Value = 2947
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8841
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -6311

$ws.Range("H136").Value = 6007796.5
$ws.Range("I136").Value = 23814222
$ws.Range("K136").Value = 71442666
$ws.Range("M136").Value = -71440116

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I34").Value = 1300
$ws.Range("J34").Value = 465.66666
$ws.Range("K34").Value = 3900
$ws.Range("L34").Value = 1396.99998
$ws.Range("M34").Value = -3816
$ws.Range("N34").Value = -1564.99998

$ws.Range("H55").Value = 77272870
$ws.Range("I55").Value = 84000160
$ws.Range("K55").Value = 252000480
$ws.Range("M55").Value = -252000303

$ws.Range("H139").Value = 7329.55
$ws.Range("I139").Value = 3968.5386
$ws.Range("K139").Value = 11905.6158
$ws.Range("M139").Value = -6765.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4186.8335
$ws.Range("I80").Value = 4133
$ws.Range("K80").Value = 4133
$ws.Range("M80").Value = -3135

$ws.Range("H83").Value = 4186.8335
$ws.Range("I83").Value = 4133
$ws.Range("K83").Value = 20665
$ws.Range("M83").Value = -15673

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 852.94116
$ws.Range("I32").Value = 5750
$ws.Range("K32").Value = 5750
$ws.Range("M32").Value = -5433

$ws.Range("H68").Value = 1683.625
$ws.Range("I68").Value = 1575.8
$ws.Range("J68").Value = 1863.3334
$ws.Range("K68").Value = 1575.8
$ws.Range("L68").Value = 1863.3334
$ws.Range("M68").Value = -826.8
$ws.Range("N68").Value = -3361.3334

$ws.Range("H71").Value = 1683.625
$ws.Range("I71").Value = 1575.8
$ws.Range("J71").Value = 1863.3334
$ws.Range("K71").Value = 7879
$ws.Range("L71").Value = 9316.666999999999
$ws.Range("M71").Value = -4135
$ws.Range("N71").Value = -16804.667

$ws.Range("H82").Value = 4162.6665
$ws.Range("I82").Value = 2489
$ws.Range("J82").Value = 4999.5
$ws.Range("K82").Value = 2489
$ws.Range("L82").Value = 4999.5
$ws.Range("M82").Value = -2128
$ws.Range("N82").Value = -5721.5

$ws.Range("H85").Value = 4162.6665
$ws.Range("I85").Value = 2489
$ws.Range("J85").Value = 4999.5
$ws.Range("K85").Value = 2489
$ws.Range("L85").Value = 4999.5
$ws.Range("M85").Value = -1241
$ws.Range("N85").Value = -7495.5

$ws.Range("H122").Value = 6995.1113
$ws.Range("I122").Value = 8199
$ws.Range("K122").Value = 24597
$ws.Range("M122").Value = -22147

$ws.Range("H132").Value = 1590939.5
$ws.Range("I132").Value = 3033739.5
$ws.Range("K132").Value = 9101218.5
$ws.Range("M132").Value = -9098688.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4903929
$ws.Range("I132").Value = 5748845
$ws.Range("J132").Value = 3417.6
$ws.Range("K132").Value = 17246535
$ws.Range("L132").Value = 10252.8
$ws.Range("M132").Value = -17244005
$ws.Range("N132").Value = -15312.8

Write-Output "Applied all Brynhildr_Profits.xlsx updates"
